$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting existing rows 35-42 down to 36-43
$ws.Rows.Item(35).Insert()

# Populate the new row 35 with its data
$ws.Cells.Item(35, 1).Value = 5
$ws.Cells.Item(35, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(35, 3).Value = "Maule"
$ws.Cells.Item(35, 4).Value = 44504
$ws.Cells.Item(35, 4).NumberFormat = $ws.Cells.Item(36, 4).NumberFormat
$ws.Cells.Item(35, 5).Value = 7
$ws.Cells.Item(35, 6).Value = 300000000
$ws.Cells.Item(35, 7).Value = "Espárragos"
$ws.Cells.Item(35, 8).Value = "Verde"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 5000
$ws.Cells.Item(35, 11).Value = 800
$ws.Cells.Item(35, 12).Value = 800
$ws.Cells.Item(35, 13).Value = 800
$ws.Cells.Item(35, 14).Value = "$/kilo"
$ws.Cells.Item(35, 15).Value = "Provincia de Linares"
$ws.Cells.Item(35, 16).Value = 800
$ws.Cells.Item(35, 17).Value = 1
$ws.Cells.Item(35, 18).Value = "Hortaliza"
